# The published sheet used to carry a small "title block" in the first
# three rows of Sheet1 (a free-form label/date row, a depth value row, and
# a blank spacer row) above the real column headers. Importing this sheet
# elsewhere choked on that extra block, so remove it: delete rows 1-3 so
# the real header row becomes row 1 and everything below shifts up to
# match.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1:C3").EntireRow.Delete() | Out-Null

# Restore the active selection on the sheet to where the user left off.
$ws.Range("C21").Select() | Out-Null
